{"js": "// Replace each two-digit-by-two-digit multiplication expression in the\n// document's table cells with its new value, per the commit's diff.\n// The (old -> new) pairs below are applied as exact, case-sensitive\n// whole-document text substitutions; each old expression occurs exactly\n// once in the document, so a direct search+replace is unambiguous.\nconst replacements = [\n  [\"87\u00d741=\", \"62\u00d788=\"],\n  [\"94\u00d795=\", \"56\u00d756=\"],\n  [\"97\u00d753=\", \"92\u00d787=\"],\n  [\"40\u00d716=\", \"41\u00d756=\"],\n  [\"91\u00d746=\", \"38\u00d733=\"],\n  [\"78\u00d758=\", \"23\u00d720=\"],\n  [\"21\u00d775=\", \"93\u00d753=\"],\n  [\"32\u00d744=\", \"47\u00d789=\"],\n  [\"69\u00d725=\", \"99\u00d768=\"],\n  [\"54\u00d776=\", \"97\u00d729=\"],\n  [\"96\u00d714=\", \"97\u00d733=\"],\n  [\"73\u00d737=\", \"50\u00d790=\"],\n  [\"33\u00d799=\", \"28\u00d744=\"],\n  [\"42\u00d778=\", \"45\u00d725=\"],\n  [\"76\u00d771=\", \"85\u00d719=\"],\n  [\"25\u00d788=\", \"21\u00d745=\"],\n  [\"96\u00d733=\", \"27\u00d771=\"],\n  [\"73\u00d724=\", \"60\u00d751=\"],\n  [\"52\u00d773=\", \"44\u00d750=\"],\n  [\"85\u00d765=\", \"33\u00d721=\"],\n  [\"88\u00d768=\", \"18\u00d788=\"],\n  [\"68\u00d721=\", \"74\u00d715=\"],\n  [\"36\u00d715=\", \"14\u00d777=\"],\n  [\"41\u00d797=\", \"16\u00d789=\"],\n  [\"34\u00d731=\", \"14\u00d782=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-by-two-digit multiplication expression in the\n# document's table cells with its new value, per the commit's diff.\n# Each old expression occurs exactly once in the document, so a\n# Find/Replace over the whole document Range is unambiguous.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"87\u00d741=\", \"62\u00d788=\"),\n    @(\"94\u00d795=\", \"56\u00d756=\"),\n    @(\"97\u00d753=\", \"92\u00d787=\"),\n    @(\"40\u00d716=\", \"41\u00d756=\"),\n    @(\"91\u00d746=\", \"38\u00d733=\"),\n    @(\"78\u00d758=\", \"23\u00d720=\"),\n    @(\"21\u00d775=\", \"93\u00d753=\"),\n    @(\"32\u00d744=\", \"47\u00d789=\"),\n    @(\"69\u00d725=\", \"99\u00d768=\"),\n    @(\"54\u00d776=\", \"97\u00d729=\"),\n    @(\"96\u00d714=\", \"97\u00d733=\"),\n    @(\"73\u00d737=\", \"50\u00d790=\"),\n    @(\"33\u00d799=\", \"28\u00d744=\"),\n    @(\"42\u00d778=\", \"45\u00d725=\"),\n    @(\"76\u00d771=\", \"85\u00d719=\"),\n    @(\"25\u00d788=\", \"21\u00d745=\"),\n    @(\"96\u00d733=\", \"27\u00d771=\"),\n    @(\"73\u00d724=\", \"60\u00d751=\"),\n    @(\"52\u00d773=\", \"44\u00d750=\"),\n    @(\"85\u00d765=\", \"33\u00d721=\"),\n    @(\"88\u00d768=\", \"18\u00d788=\"),\n    @(\"68\u00d721=\", \"74\u00d715=\"),\n    @(\"36\u00d715=\", \"14\u00d777=\"),\n    @(\"41\u00d797=\", \"16\u00d789=\"),\n    @(\"34\u00d731=\", \"14\u00d782=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
